$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 25 de Marzo de 2020 a las 00:46"

# Row 6: Estados Unidos
$ws.Range("B6").Value = 53358
$ws.Range("C6").Value = 9624
$ws.Range("D6").Value = 370
$ws.Range("E6").Value = 52293
$ws.Range("F6").Value = 1175
$ws.Range("G6").Value = 142
$ws.Range("H6").Value = 695

# Row 7: España
$ws.Range("B7").Value = 42058
$ws.Range("C7").Value = 6922
$ws.Range("D7").Value = 3794
$ws.Range("E7").Value = 35273
$ws.Range("F7").Value = 2636
$ws.Range("G7").Value = 680
$ws.Range("H7").Value = 2991

# Row 31: Ecuador
$ws.Range("B31").Value = 1082
$ws.Range("C31").Value = 101
$ws.Range("D31").Value = 3
$ws.Range("E31").Value = 1052
$ws.Range("F31").Value = 2
$ws.Range("G31").Value = 9
$ws.Range("H31").Value = 27

# Row 53: Argentina
$ws.Range("A53").Value = "Argentina"
$ws.Range("B53").Value = 387
$ws.Range("C53").Value = 86
$ws.Range("D53").Value = 52
$ws.Range("E53").Value = 329
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 2
$ws.Range("H53").Value = 6

# Row 54: Hong Kong
$ws.Range("A54").Value = "Hong Kong"
$ws.Range("B54").Value = 386
$ws.Range("C54").Value = 29
$ws.Range("D54").Value = 102
$ws.Range("E54").Value = 280
$ws.Range("F54").Value = 4
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 4

# Row 55: Croacia
$ws.Range("A55").Value = "Croacia"
$ws.Range("B55").Value = 382
$ws.Range("C55").Value = 67
$ws.Range("D55").Value = 5
$ws.Range("E55").Value = 376
$ws.Range("F55").Value = 6
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 1

# Row 56: Colombia
$ws.Range("A56").Value = "Colombia"
$ws.Range("B56").Value = 378
$ws.Range("C56").Value = 101
$ws.Range("D56").Value = 6
$ws.Range("E56").Value = 369
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 3

# Row 57: Estonia
$ws.Range("A57").Value = "Estonia"
$ws.Range("B57").Value = 369
$ws.Range("C57").Value = 17
$ws.Range("D57").Value = 7
$ws.Range("E57").Value = 362
$ws.Range("F57").Value = 4
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 0

# Row 58: Mexico
$ws.Range("A58").Value = "Mexico"
$ws.Range("B58").Value = 367
$ws.Range("C58").Value = 51
$ws.Range("D58").Value = 4
$ws.Range("E58").Value = 359
$ws.Range("F58").Value = 1
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 4

# Row 59: Panama
$ws.Range("A59").Value = "Panama"
$ws.Range("B59").Value = 345
$ws.Range("C59").Value = 0
$ws.Range("D59").Value = 1
$ws.Range("E59").Value = 338
$ws.Range("F59").Value = 33
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 6

# Row 60: Libano
$ws.Range("A60").Value = "Libano"
$ws.Range("B60").Value = 318
$ws.Range("C60").Value = 51
$ws.Range("D60").Value = 8
$ws.Range("E60").Value = 306
$ws.Range("F60").Value = 4
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 4

# Row 61: Irak
$ws.Range("A61").Value = "Irak"
$ws.Range("B61").Value = 316
$ws.Range("C61").Value = 50
$ws.Range("D61").Value = 75
$ws.Range("E61").Value = 214
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 4
$ws.Range("H61").Value = 27

# Row 62: Republica Dominicana
$ws.Range("A62").Value = "Republica Dominicana"
$ws.Range("B62").Value = 312
$ws.Range("C62").Value = 67
$ws.Range("D62").Value = 3
$ws.Range("E62").Value = 303
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 3
$ws.Range("H62").Value = 6

# Row 63: Serbia
$ws.Range("A63").Value = "Serbia"
$ws.Range("B63").Value = 303
$ws.Range("C63").Value = 54
$ws.Range("D63").Value = 15
$ws.Range("E63").Value = 285
$ws.Range("F63").Value = 21
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 3

# Row 102: Guadalupe
$ws.Range("A102").Value = "Guadalupe"
$ws.Range("B102").Value = 73
$ws.Range("C102").Value = 11
$ws.Range("D102").Value = 0
$ws.Range("E102").Value = 72
$ws.Range("F102").Value = 4
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 1

# Row 103: Costa de Marfil
$ws.Range("A103").Value = "Costa de Marfil"
$ws.Range("B103").Value = 73
$ws.Range("C103").Value = 48
$ws.Range("D103").Value = 2
$ws.Range("E103").Value = 71
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 0

# Row 104: Kazajistan
$ws.Range("A104").Value = "Kazajistan"
$ws.Range("B104").Value = 72
$ws.Range("C104").Value = 10
$ws.Range("D104").Value = 0
$ws.Range("E104").Value = 72
$ws.Range("F104").Value = 0
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 0

# Row 105: Georgia
$ws.Range("A105").Value = "Georgia"
$ws.Range("B105").Value = 70
$ws.Range("C105").Value = 9
$ws.Range("D105").Value = 9
$ws.Range("E105").Value = 61
$ws.Range("F105").Value = 1
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 0

# Row 106: Camerun
$ws.Range("A106").Value = "Camerun"
$ws.Range("B106").Value = 66
$ws.Range("C106").Value = 10
$ws.Range("D106").Value = 2
$ws.Range("E106").Value = 64
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 0

# Row 137: Islas Virgenes de los Estados Unidos
$ws.Range("A137").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("B137").Value = 17
$ws.Range("C137").Value = 0
$ws.Range("D137").Value = 0
$ws.Range("E137").Value = 17
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 0

# Row 138: Madagascar
$ws.Range("A138").Value = "Madagascar"
$ws.Range("B138").Value = 17
$ws.Range("C138").Value = 5
$ws.Range("D138").Value = 0
$ws.Range("E138").Value = 17
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 0

# Row 141: Tanzania
$ws.Range("A141").Value = "Tanzania"
$ws.Range("B141").Value = 12
$ws.Range("C141").Value = 0
$ws.Range("D141").Value = 0
$ws.Range("E141").Value = 12
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 0

# Row 142: Etiopia
$ws.Range("A142").Value = "Etiopia"
$ws.Range("B142").Value = 12
$ws.Range("C142").Value = 1
$ws.Range("D142").Value = 0
$ws.Range("E142").Value = 12
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 0

# Row 144: Nueva Caledonia
$ws.Range("A144").Value = "Nueva Caledonia"
$ws.Range("B144").Value = 10
$ws.Range("C144").Value = 2
$ws.Range("D144").Value = 0
$ws.Range("E144").Value = 10
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 0

# Row 145: Mongolia
$ws.Range("A145").Value = "Mongolia"
$ws.Range("B145").Value = 10
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 0
$ws.Range("E145").Value = 10
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 0

# Row 146: Uganda
$ws.Range("A146").Value = "Uganda"
$ws.Range("B146").Value = 9
$ws.Range("C146").Value = 0
$ws.Range("D146").Value = 0
$ws.Range("E146").Value = 9
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 0

# Row 147: Guinea Ecuatorial
$ws.Range("A147").Value = "Guinea Ecuatorial"
$ws.Range("B147").Value = 9
$ws.Range("C147").Value = 0
$ws.Range("D147").Value = 0
$ws.Range("E147").Value = 9
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 0

# Row 149: Haiti
$ws.Range("A149").Value = "Haiti"
$ws.Range("B149").Value = 7
$ws.Range("C149").Value = 1
$ws.Range("D149").Value = 0
$ws.Range("E149").Value = 7
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 0

# Row 150: Seychelles
$ws.Range("A150").Value = "Seychelles"
$ws.Range("B150").Value = 7
$ws.Range("C150").Value = 0
$ws.Range("D150").Value = 0
$ws.Range("E150").Value = 7
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 0

# Row 151: Surinam
$ws.Range("A151").Value = "Surinam"
$ws.Range("B151").Value = 7
$ws.Range("C151").Value = 2
$ws.Range("D151").Value = 0
$ws.Range("E151").Value = 7
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 0

# Row 155: Curazao
$ws.Range("A155").Value = "Curazao"
$ws.Range("B155").Value = 6
$ws.Range("C155").Value = 2
$ws.Range("D155").Value = 0
$ws.Range("E155").Value = 5
$ws.Range("F155").Value = 0
$ws.Range("G155").Value = 0
$ws.Range("H155").Value = 1

# Row 156: Gabon
$ws.Range("A156").Value = "Gabon"
$ws.Range("B156").Value = 6
$ws.Range("C156").Value = 0
$ws.Range("D156").Value = 0
$ws.Range("E156").Value = 5
$ws.Range("F156").Value = 0
$ws.Range("G156").Value = 0
$ws.Range("H156").Value = 1

# Row 159: Bahamas
$ws.Range("A159").Value = "Bahamas"
$ws.Range("B159").Value = 5
$ws.Range("C159").Value = 1
$ws.Range("D159").Value = 1
$ws.Range("E159").Value = 4
$ws.Range("F159").Value = 0
$ws.Range("G159").Value = 0
$ws.Range("H159").Value = 0

# Row 160: Guyana
$ws.Range("A160").Value = "Guyana"
$ws.Range("B160").Value = 5
$ws.Range("C160").Value = 0
$ws.Range("D160").Value = 0
$ws.Range("E160").Value = 4
$ws.Range("F160").Value = 0
$ws.Range("G160").Value = 0
$ws.Range("H160").Value = 1

# Row 165: Congo
$ws.Range("A165").Value = "Congo"
$ws.Range("B165").Value = 4
$ws.Range("C165").Value = 0
$ws.Range("D165").Value = 0
$ws.Range("E165").Value = 4
$ws.Range("F165").Value = 0
$ws.Range("G165").Value = 0
$ws.Range("H165").Value = 0

# Row 166: Fiyi
$ws.Range("A166").Value = "Fiyi"
$ws.Range("B166").Value = 4
$ws.Range("C166").Value = 1
$ws.Range("D166").Value = 0
$ws.Range("E166").Value = 4
$ws.Range("F166").Value = 0
$ws.Range("G166").Value = 0
$ws.Range("H166").Value = 0

# Row 167: Zambia
$ws.Range("A167").Value = "Zambia"
$ws.Range("B167").Value = 3
$ws.Range("C167").Value = 0
$ws.Range("D167").Value = 0
$ws.Range("E167").Value = 3
$ws.Range("F167").Value = 0
$ws.Range("G167").Value = 0
$ws.Range("H167").Value = 0

# Row 170: Santa Lucia
$ws.Range("A170").Value = "Santa Lucia"
$ws.Range("B170").Value = 3
$ws.Range("C170").Value = 0
$ws.Range("D170").Value = 0
$ws.Range("E170").Value = 3
$ws.Range("F170").Value = 0
$ws.Range("G170").Value = 0
$ws.Range("H170").Value = 0

# Row 171: Birmania
$ws.Range("A171").Value = "Birmania"
$ws.Range("B171").Value = 3
$ws.Range("C171").Value = 1
$ws.Range("D171").Value = 0
$ws.Range("E171").Value = 3
$ws.Range("F171").Value = 0
$ws.Range("G171").Value = 0
$ws.Range("H171").Value = 0

# Row 172: San Bartolome
$ws.Range("A172").Value = "San Bartolome"
$ws.Range("B172").Value = 3
$ws.Range("C172").Value = 0
$ws.Range("D172").Value = 0
$ws.Range("E172").Value = 3
$ws.Range("F172").Value = 0
$ws.Range("G172").Value = 0
$ws.Range("H172").Value = 0

# Row 173: Niger
$ws.Range("A173").Value = "Niger"
$ws.Range("B173").Value = 3
$ws.Range("C173").Value = 0
$ws.Range("D173").Value = 0
$ws.Range("E173").Value = 3
$ws.Range("F173").Value = 0
$ws.Range("G173").Value = 0
$ws.Range("H173").Value = 0

# Row 174: Republica de Africa Central
$ws.Range("A174").Value = "Republica de Africa Central"
$ws.Range("B174").Value = 3
$ws.Range("C174").Value = 0
$ws.Range("D174").Value = 0
$ws.Range("E174").Value = 3
$ws.Range("F174").Value = 0
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 0

# Row 175: Angola
$ws.Range("A175").Value = "Angola"
$ws.Range("B175").Value = 3
$ws.Range("C175").Value = 0
$ws.Range("D175").Value = 0
$ws.Range("E175").Value = 3
$ws.Range("F175").Value = 0
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = 0

# Row 176: Liberia
$ws.Range("A176").Value = "Liberia"
$ws.Range("B176").Value = 3
$ws.Range("C176").Value = 0
$ws.Range("D176").Value = 0
$ws.Range("E176").Value = 3
$ws.Range("F176").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 0

# Row 178: Republica del Chad
$ws.Range("A178").Value = "Republica del Chad"
$ws.Range("B178").Value = 3
$ws.Range("C178").Value = 1
$ws.Range("D178").Value = 0
$ws.Range("E178").Value = 3
$ws.Range("F178").Value = 0
$ws.Range("G178").Value = 0
$ws.Range("H178").Value = 0

# Row 179: Gambia
$ws.Range("A179").Value = "Gambia"
$ws.Range("B179").Value = 3
$ws.Range("C179").Value = 1
$ws.Range("D179").Value = 0
$ws.Range("E179").Value = 2
$ws.Range("F179").Value = 0
$ws.Range("G179").Value = 0
$ws.Range("H179").Value = 1

# Row 181: Zimbabue
$ws.Range("A181").Value = "Zimbabue"
$ws.Range("B181").Value = 3
$ws.Range("C181").Value = 0
$ws.Range("D181").Value = 0
$ws.Range("E181").Value = 2
$ws.Range("F181").Value = 0
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 1

# Row 183: Dominica
$ws.Range("A183").Value = "Dominica"
$ws.Range("B183").Value = 2
$ws.Range("C183").Value = 0
$ws.Range("D183").Value = 0
$ws.Range("E183").Value = 2
$ws.Range("F183").Value = 0
$ws.Range("G183").Value = 0
$ws.Range("H183").Value = 0

# Row 184: Nicaragua
$ws.Range("A184").Value = "Nicaragua"
$ws.Range("B184").Value = 2
$ws.Range("C184").Value = 0
$ws.Range("D184").Value = 0
$ws.Range("E184").Value = 2
$ws.Range("F184").Value = 0
$ws.Range("G184").Value = 0
$ws.Range("H184").Value = 0

# Row 185: Mauritania
$ws.Range("A185").Value = "Mauritania"
$ws.Range("B185").Value = 2
$ws.Range("C185").Value = 0
$ws.Range("D185").Value = 0
$ws.Range("E185").Value = 2
$ws.Range("F185").Value = 0
$ws.Range("G185").Value = 0
$ws.Range("H185").Value = 0

# Row 186: Laos
$ws.Range("A186").Value = "Laos"
$ws.Range("B186").Value = 2
$ws.Range("C186").Value = 2
$ws.Range("D186").Value = 0
$ws.Range("E186").Value = 2
$ws.Range("F186").Value = 0
$ws.Range("G186").Value = 0
$ws.Range("H186").Value = 0

# Row 187: Butan
$ws.Range("A187").Value = "Butan"
$ws.Range("B187").Value = 2
$ws.Range("C187").Value = 0
$ws.Range("D187").Value = 0
$ws.Range("E187").Value = 2
$ws.Range("F187").Value = 0
$ws.Range("G187").Value = 0
$ws.Range("H187").Value = 0

# Row 188: San Martin (Parte Holandesa)
$ws.Range("A188").Value = "San Martin (Parte Holandesa)"
$ws.Range("B188").Value = 2
$ws.Range("C188").Value = 0
$ws.Range("D188").Value = 0
$ws.Range("E188").Value = 2
$ws.Range("F188").Value = 0
$ws.Range("G188").Value = 0
$ws.Range("H188").Value = 0

# Row 190: Somalia
$ws.Range("A190").Value = "Somalia"
$ws.Range("B190").Value = 1
$ws.Range("C190").Value = 0
$ws.Range("D190").Value = 0
$ws.Range("E190").Value = 1
$ws.Range("F190").Value = 0
$ws.Range("G190").Value = 0
$ws.Range("H190").Value = 0

# Row 191: Timor Oriental
$ws.Range("A191").Value = "Timor Oriental"
$ws.Range("B191").Value = 1
$ws.Range("C191").Value = 0
$ws.Range("D191").Value = 0
$ws.Range("E191").Value = 1
$ws.Range("F191").Value = 0
$ws.Range("G191").Value = 0
$ws.Range("H191").Value = 0

# Row 192: Libia
$ws.Range("A192").Value = "Libia"
$ws.Range("B192").Value = 1
$ws.Range("C192").Value = 1
$ws.Range("D192").Value = 0
$ws.Range("E192").Value = 1
$ws.Range("F192").Value = 0
$ws.Range("G192").Value = 0
$ws.Range("H192").Value = 0

# Row 194: Montserrat
$ws.Range("A194").Value = "Montserrat"
$ws.Range("B194").Value = 1
$ws.Range("C194").Value = 0
$ws.Range("D194").Value = 0
$ws.Range("E194").Value = 1
$ws.Range("F194").Value = 0
$ws.Range("G194").Value = 0
$ws.Range("H194").Value = 0

# Row 195: Belice
$ws.Range("A195").Value = "Belice"
$ws.Range("B195").Value = 1
$ws.Range("C195").Value = 0
$ws.Range("D195").Value = 0
$ws.Range("E195").Value = 1
$ws.Range("F195").Value = 0
$ws.Range("G195").Value = 0
$ws.Range("H195").Value = 0

# Row 196: San Vicente y las Granadinas
$ws.Range("A196").Value = "San Vicente y las Granadinas"
$ws.Range("B196").Value = 1
$ws.Range("C196").Value = 0
$ws.Range("D196").Value = 0
$ws.Range("E196").Value = 1
$ws.Range("F196").Value = 0
$ws.Range("G196").Value = 0
$ws.Range("H196").Value = 0

# Row 197: Siria
$ws.Range("A197").Value = "Siria"
$ws.Range("B197").Value = 1
$ws.Range("C197").Value = 0
$ws.Range("D197").Value = 0
$ws.Range("E197").Value = 1
$ws.Range("F197").Value = 0
$ws.Range("G197").Value = 0
$ws.Range("H197").Value = 0

# Row 198: Papua Nueva Guinea
$ws.Range("A198").Value = "Papua Nueva Guinea"
$ws.Range("B198").Value = 1
$ws.Range("C198").Value = 0
$ws.Range("D198").Value = 0
$ws.Range("E198").Value = 1
$ws.Range("F198").Value = 0
$ws.Range("G198").Value = 0
$ws.Range("H198").Value = 0

# Row 199: Granada
$ws.Range("A199").Value = "Granada"
$ws.Range("B199").Value = 1
$ws.Range("C199").Value = 0
$ws.Range("D199").Value = 0
$ws.Range("E199").Value = 1
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 0

# Row 200: Islas Turcas y Caicos
$ws.Range("A200").Value = "Islas Turcas y Caicos"
$ws.Range("B200").Value = 1
$ws.Range("C200").Value = 0
$ws.Range("D200").Value = 0
$ws.Range("E200").Value = 1
$ws.Range("F200").Value = 0
$ws.Range("G200").Value = 0
$ws.Range("H200").Value = 0

